# Refresh the cryptos price/volume snapshot in-place (GitHub Actions style update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.394.66'
$ws.Range('E2').Value = '  +0.03%  '

$ws.Range('D3').Value = '1.799.56'
$ws.Range('E3').Value = '  -0.85%  '

$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').Value = '''227.27'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.34%  '

$ws.Range('D6').Value = '''0.573'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.82%  '

$ws.Range('E7').Value = '  +0.21%  '

$ws.Range('D8').Value = '''36.16'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +6.96%  '

$ws.Range('E9').Value = '  +0.33%  '

$ws.Range('D10').Value = '''0.0690'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.25%  '

$ws.Range('D11').Value = '''0.0964'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.77%  '

$ws.Range('D12').Value = '2.060.28'
$ws.Range('E12').Value = '  -0.31%  '

$ws.Range('D13').Value = '''11.54'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.69%  '

$ws.Range('D14').Value = '''0.644'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.47%  '

$ws.Range('D15').Value = '1.769.34'
$ws.Range('E15').Value = '  -1.99%  '

$ws.Range('E16').Value = '  +4.53%  '

$ws.Range('D17').Value = '34.364.82'
$ws.Range('E17').Value = '  -0.02%  '

$ws.Range('D18').Value = '''68.89'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.26%  '

$ws.Range('D19').Value = '''244.44'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.41%  '

$ws.Range('D20').Value = '0.0₃0790'
$ws.Range('E20').Value = '  -1.06%  '

$ws.Range('E21').Value = '  +1.74%  '

$ws.Range('E22').Value = '  -0.08%  '

$ws.Range('D23').Value = '''4.17'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.40%  '

$ws.Range('D24').Value = '''172.66'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.06%  '

$ws.Range('E25').Value = '  +3.46%  '

$ws.Range('D26').Value = '''7.96'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +8.15%  '

$ws.Range('D27').Value = '''16.79'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.11%  '

$ws.Range('E28').Value = '  +1.27%  '

$ws.Range('E29').Value = '  +0.09%  '

$ws.Range('D30').Value = '''4.01'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.33%  '

$ws.Range('E31').Value = '  -0.23%  '

$ws.Range('D32').Value = '''3.83'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.35%  '

$ws.Range('E33').Value = '  -0.39%  '

$ws.Range('E34').Value = '  -1.89%  '

$ws.Range('D35').Value = '1.396.61'
$ws.Range('E35').Value = '  -1.12%  '

$ws.Range('D36').Value = '''0.671'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.95%  '

$ws.Range('D37').Value = '''2.44'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -6.77%  '

$ws.Range('E38').Value = '  -0.56%  '

$ws.Range('D39').Value = '''0.0190'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.60%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '''82.36'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.92%  '

$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '''0.957'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.43%  '

$ws.Range('E42').Value = '  -1.05%  '

$ws.Range('D43').Value = '''2.42'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.48%  '

$ws.Range('E44').Value = '  +6.49%  '

$ws.Range('D45').Value = '''13.30'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.86%  '

$ws.Range('D46').Value = '''6.02'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.72%  '

$ws.Range('D47').Value = '''0.0503'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.81%  '

$ws.Range('D48').Value = '1.961.35'
$ws.Range('E48').Value = '  +0.21%  '

$ws.Range('D49').Value = '''104.22'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.38%  '

$ws.Range('E50').Value = '  +0.14%  '

$ws.Range('E51').Value = '  +0.92%  '
